$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.905.96'
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').Value = '2.601.00'
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '523.29'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +3.01%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '154.87'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.77%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.588'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.99%  '
$ws.Range('E9').Value = '  +2.39%  '
$ws.Range('E10').Value = '  +1.35%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.348'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.18%  '
$ws.Range('E12').Value = '  +1.01%  '
$ws.Range('D13').Value = '3.058.81'
$ws.Range('E13').Value = '  +0.61%  '
$ws.Range('D14').Value = '60.933.29'
$ws.Range('E14').Value = '  +0.61%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '21.70'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.12%  '
$ws.Range('E16').Value = '  +0.83%  '
$ws.Range('D17').Value = '2.604.03'
$ws.Range('E17').Value = '  +0.70%  '
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '355.26'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.52%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '10.59'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.33%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.22'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.91%  '
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '61.05'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.67%  '
$ws.Range('E24').Value = '  +1.42%  '
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('D26').Value = '2.717.12'
$ws.Range('E26').Value = '  +0.63%  '
$ws.Range('E27').Value = '  +0.23%  '
$ws.Range('E28').Value = '  +0.07%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.27'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +9.62%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '19.43'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('E33').Value = '  +2.88%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '148.89'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -2.80%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.21'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +5.42%  '
$ws.Range('E36').Value = '  +1.15%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.912'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +7.34%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.895'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +5.48%  '
$ws.Range('E39').Value = '  +1.42%  '
$ws.Range('E40').Value = '  +1.53%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '291.07'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.57%  '
$ws.Range('E43').Value = '  +2.07%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.623'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('E45').Value = '  +0.44%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.999'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.20%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '19.58'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.12%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '4.93'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.89%  '
$ws.Range('E49').Value = '  +2.01%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '19.27'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +9.24%  '
